# NeapWorkflow GET MarineTerrestrial Pre-1750 250m workbook update
# Commit: "Updated all workflow files with new path to NVIS Extant and Pre1750 files"
#
# The NVIS pre-1750 intermediate raster was reprocessed on 2024-08-01 (was 2024-07-30),
# so the RawDataPath cell for the "Terrestrial-Pre-IUCNGET" row needs to point at the
# new file. That cell lives at Data!B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B4").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_PRE1750_IUCNGET_DK_20240801.tif"

# Reflect the author's saved UI state: the active sheet/cell moved from A4 to B4
# (the long RawDataPath column) with the view scrolled one column to the right.
$ws.Activate()
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
